$d = $word.ActiveDocument

# 1. Update the "finish by" date clause
$d.Content.Find.Execute(
    "see it finished by December 14th 2016.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "finish it by January 8th 2018.", 2
)

# 2. Remove the trailing "online meetings" sentence
$d.Content.Find.Execute(
    " This time will be supported by any online meetings taken by me and my partner in which we work on and discuss the project as needed throughout the week.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " ", 2
)

# 3. Update the gameplay description: one map, two characters hitting a ball
$d.Content.Find.Execute(
    "two playable maps which graphical design will be developed by my partner and will consist of a ball which moves around the screen and collides with objects and game boundaries causing different events to be initiated.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "one playable map which graphical design will be developed by my partner and will consist of two playable characters which have to hit a ball around the screen in order to hit the other playable character.", 2
)
